$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Correct the C3 value (it was entered as a whole number "18" instead of a
# fraction "0.18").
$ws.Range("C3").Value = 0.18

# Move the active selection back to F2 (top of sheet) instead of leaving it
# further down at F22.
$ws.Range("F2").Select()
